# data : case 1
# Apply the column-width and cell-value changes described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes -------------------------------------------------
# Target widths (raw OOXML "width" units): D/H/I -> 2.140625, N -> 4.7109375
# Excel's ColumnWidth property is quantized to the workbook's Normal-style
# pixel grid, so we pick the ColumnWidth input that lands on the closest
# achievable width to each target.
$ws.Columns.Item(4).ColumnWidth = 1.33    # column D: 3.140625 -> ~2.140625
$ws.Columns.Item(8).ColumnWidth = 1.33    # column H: 3.140625 -> ~2.140625
$ws.Columns.Item(9).ColumnWidth = 1.33    # column I: 3.140625 -> ~2.140625
$ws.Columns.Item(14).ColumnWidth = 3.83   # column N: 5.7109375 -> ~4.7109375

# --- Cell value changes (row 1) -------------------------------------------
$ws.Range("C1").Value = 12
$ws.Range("D1").Value = 7
$ws.Range("E1").Value = 22
$ws.Range("F1").Value = 18
$ws.Range("G1").Value = 15
$ws.Range("H1").Value = 3
$ws.Range("I1").Value = 5
$ws.Range("J1").Value = 32
$ws.Range("K1").Value = 0.052000000000000005
$ws.Range("L1").Value = 0.087
$ws.Range("M1").Value = 0.096
$ws.Range("N1").Value = 0.06999999999999999
